$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.616.31"
$ws.Range("E2").Value = "  +4.56%  "

$ws.Range("D3").Value = "3.486.17"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'591.87"
$ws.Range("E5").Value = "  +3.99%  "

$ws.Range("D6").Value = "'169.30"
$ws.Range("E6").Value = "  +5.35%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "3.482.90"
$ws.Range("E8").Value = "  +2.73%  "

$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  +8.70%  "

$ws.Range("D10").Value = "'7.32"
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("E11").Value = "  +7.64%  "

$ws.Range("E12").Value = "  +4.77%  "

$ws.Range("D13").Value = "4.090.75"
$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "'28.13"
$ws.Range("E15").Value = "  +5.21%  "

$ws.Range("E16").Value = "  +4.25%  "

$ws.Range("D17").Value = "66.628.16"
$ws.Range("E17").Value = "  +4.45%  "

$ws.Range("D18").Value = "3.476.98"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("E19").Value = "  +3.91%  "

$ws.Range("D20").Value = "'14.03"
$ws.Range("E20").Value = "  +4.39%  "

$ws.Range("D21").Value = "'391.88"
$ws.Range("E21").Value = "  +4.55%  "

$ws.Range("E22").Value = "  +2.28%  "

$ws.Range("D23").Value = "'73.05"
$ws.Range("E23").Value = "  +4.52%  "

$ws.Range("E25").Value = "  +4.96%  "

$ws.Range("E26").Value = "  +7.53%  "

$ws.Range("D27").Value = "'10.29"
$ws.Range("E27").Value = "  +8.86%  "

$ws.Range("E28").Value = "  +1.59%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +5.13%  "

$ws.Range("D31").Value = "'1.46"
$ws.Range("E31").Value = "  +5.93%  "

$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +3.76%  "

$ws.Range("E33").Value = "  +3.75%  "

$ws.Range("D34").Value = "'7.43"
$ws.Range("E34").Value = "  +6.64%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("E36").Value = "  +10.37%  "

$ws.Range("D37").Value = "'161.54"
$ws.Range("E37").Value = "  +1.24%  "

$ws.Range("D38").Value = "'0.903"
$ws.Range("E38").Value = "  +5.36%  "

$ws.Range("E39").Value = "  +7.72%  "

$ws.Range("D40").Value = "'6.78"
$ws.Range("E40").Value = "  +6.48%  "

$ws.Range("E41").Value = "  +3.92%  "

$ws.Range("D42").Value = "'26.63"
$ws.Range("E42").Value = "  +4.09%  "

$ws.Range("D43").Value = "'4.64"
$ws.Range("E43").Value = "  +7.30%  "

$ws.Range("D44").Value = "'26.79"
$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("D45").Value = "'43.18"
$ws.Range("E45").Value = "  +1.39%  "

$ws.Range("D46").Value = "2.762.70"
$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("E47").Value = "  +3.06%  "

$ws.Range("E48").Value = "  +4.76%  "

$ws.Range("D49").Value = "'346.37"
$ws.Range("E49").Value = "  +6.12%  "

$ws.Range("E50").Value = "  +5.71%  "

$ws.Range("D51").Value = "'33.90"
$ws.Range("E51").Value = "  +13.91%  "
